$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Objects of Interest table (rows 18-32): strip the folder prefixes
# from the info/hitmap/highlight/found-image filename columns.
$ooiRows = 18..32

# New unique strings get appended to the shared-string table in the order
# they are first assigned, so touch all of column D first, then I, then J,
# then K to reproduce the canonical ordering (info.txt*, hitmap.png*,
# highlight.png*, saturated.png*).
foreach ($r in $ooiRows) {
    $dCell = $ws.Cells.Item($r, 4)   # D: info_text_file
    $dCell.Value = $dCell.Value().Replace("assets/ooi info/", "")
}
foreach ($r in $ooiRows) {
    $iCell = $ws.Cells.Item($r, 9)   # I: hitmap_filename
    $iCell.Value = $iCell.Value().Replace("assets/ooi images/", "")
}
foreach ($r in $ooiRows) {
    $jCell = $ws.Cells.Item($r, 10)  # J: highlight_filename
    $jCell.Value = $jCell.Value().Replace("assets/ooi images/", "")
}
foreach ($r in $ooiRows) {
    $kCell = $ws.Cells.Item($r, 11)  # K: found_image_filename
    $kCell.Value = $kCell.Value().Replace("assets/ooi images/", "")
}

# Parrott rifle / row 18 "y" coordinate tweak
$ws.Cells.Item(18, 6).Value = 0.70486368635999996

# Shrink the table columns to fit the now-shorter filenames (best effort -
# the COM layer quantizes ColumnWidth to whole pixels, so these are chosen
# to land as close as possible to the canonical post-edit widths).
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(2).ColumnWidth = 56.833333333333336
$ws.Columns.Item(3).ColumnWidth = 75.5
$ws.Columns.Item(4).ColumnWidth = 17.666666666666668
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 13.5
$ws.Columns.Item(9).ColumnWidth = 21.5
$ws.Columns.Item(10).ColumnWidth = 23.0
$ws.Columns.Item(11).ColumnWidth = 23.666666666666668

# Update the active selection to match the new edit location
$ws.Range("E18").Select()
